$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, so Excel keeps them as strings (matching the source data which
# stores all price/volume figures as inline text) instead of auto-converting them
# to numeric cells.
$textCells = @("D5", "D6", "D10", "D13", "D18", "D19", "D20", "D21", "D26", "D30", "D31", "D34", "D36", "D40", "D41", "D42", "D45", "D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cryptocurrency price / 1h volume change values.
$ws.Range("D2").Value = "59.804.23"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.362.07"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "558.38"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "133.34"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("D10").Value = "5.63"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("E12").Value = "  -3.45%  "
$ws.Range("D13").Value = "24.23"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "2.786.72"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "59.779.12"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "2.363.44"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "11.06"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "4.46"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "320.08"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("D21").Value = "6.65"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "8.39"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").Value = "0.0₃0756"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "170.54"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "6.07"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "18.11"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D36").Value = "1.32"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").Value = "317.49"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "38.57"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "144.56"
$ws.Range("E42").Value = "  +3.00%  "
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "19.35"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").Value = "11.07"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  -1.82%  "
